# Mẫu 26 - Biên bản xác nhận tiến độ
# Remove the leftover "vnpt.SiteAddress" merge-field placeholder text that
# follows "Địa chỉ: " in the document header/info block. The run carrying
# this placeholder text is deleted entirely (not just its text), matching
# the authored edit.

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "vnpt.SiteAddress",  # FindText
    $true,                # MatchCase
    $true,                # MatchWholeWord
    $false,               # MatchWildcards
    $false,               # MatchSoundsLike
    $false,               # MatchAllWordForms
    $true,                # Forward
    1,                    # Wrap (wdFindContinue)
    $false,               # Format
    "",                   # ReplaceWith (remove the text -> empty run)
    2                     # Replace (wdReplaceAll)
)
